$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 16067.143
$ws.Range("J43").Value = 2304
$ws.Range("L43").Value = 2304
$ws.Range("N43").Value = -2442
$ws.Range("H76").Value = 7003
$ws.Range("I76").Value = 7003
$ws.Range("K76").Value = 7003
$ws.Range("M76").Value = -6688
$ws.Range("H79").Value = 7003
$ws.Range("I79").Value = 7003
$ws.Range("K79").Value = 7003
$ws.Range("M79").Value = -5911
$ws.Range("H80").Value = 800.26666
$ws.Range("I80").Value = 628.26666
$ws.Range("K80").Value = 1884.79998
$ws.Range("M80").Value = -886.79998
$ws.Range("H83").Value = 800.26666
$ws.Range("I83").Value = 628.26666
$ws.Range("K83").Value = 5654.39994
$ws.Range("M83").Value = -662.3999400000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 35000
$ws.Range("I19").Value = 35000
$ws.Range("K19").Value = 35000
$ws.Range("M19").Value = -34771
$ws.Range("H23").Value = 49500
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H31").Value = 7573.8
$ws.Range("I31").Value = 7573.8
$ws.Range("K31").Value = 7573.8
$ws.Range("M31").Value = -7279.8
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H97").Value = 795.73914
$ws.Range("I97").Value = 830.8889
$ws.Range("K97").Value = 830.8889
$ws.Range("M97").Value = -334.8889
$ws.Range("H132").Value = 3394.4666
$ws.Range("I132").Value = 2760.3
$ws.Range("K132").Value = 8280.900000000001
$ws.Range("M132").Value = -5750.900000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 400282.38
$ws.Range("I22").Value = 651.1539
$ws.Range("J22").Value = 500190.2
$ws.Range("K22").Value = 651.1539
$ws.Range("L22").Value = 500190.2
$ws.Range("M22").Value = -478.1539
$ws.Range("N22").Value = -500536.2
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("N57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("N136").Value = 0
$ws.Range("L136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 5751
$ws.Range("I2").Value = 4499.5
$ws.Range("J2").Value = 7002.5
$ws.Range("K2").Value = 4499.5
$ws.Range("L2").Value = 7002.5
$ws.Range("M2").Value = -4386.5
$ws.Range("N2").Value = -7228.5
$ws.Range("H22").Value = 281.66666
$ws.Range("H31").Value = 8380567
$ws.Range("I31").Value = 15209761
$ws.Range("J31").Value = 33774.445
$ws.Range("K31").Value = 15209761
$ws.Range("L31").Value = 33774.445
$ws.Range("M31").Value = -15209466
$ws.Range("N31").Value = -34364.445
$ws.Range("H34").Value = 8380567
$ws.Range("I34").Value = 15209761
$ws.Range("J34").Value = 33774.445
$ws.Range("K34").Value = 15209761
$ws.Range("L34").Value = 33774.445
$ws.Range("M34").Value = -15209559
$ws.Range("N34").Value = -34178.445
$ws.Range("H41").Value = 10999.454
$ws.Range("J41").Value = 10999.454
$ws.Range("L41").Value = 10999.454
$ws.Range("N41").Value = -11855.454
$ws.Range("H132").Value = 3507.1333
$ws.Range("I132").Value = 3147.4546
$ws.Range("K132").Value = 9442.363799999999
$ws.Range("M132").Value = -6912.363799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9252.5
$ws.Range("I5").Value = 468.33334
$ws.Range("K5").Value = 1405.00002
$ws.Range("M5").Value = -1293.00002
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H132").Value = 1460
$ws.Range("I132").Value = 1138.6111
$ws.Range("J132").Value = 2286.4285
$ws.Range("K132").Value = 10247.4999
$ws.Range("L132").Value = 20577.8565
$ws.Range("M132").Value = -7717.499900000001
$ws.Range("N132").Value = -25637.8565
$ws.Range("H135").Value = 9252.5
$ws.Range("I135").Value = 468.33334
$ws.Range("K135").Value = 4215.00006
$ws.Range("M135").Value = -1680.00006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1476.25
$ws.Range("I80").Value = 1120.6666
$ws.Range("J80").Value = 2543
$ws.Range("K80").Value = 1120.6666
$ws.Range("L80").Value = 2543
$ws.Range("M80").Value = -122.6666
$ws.Range("N80").Value = -4539
$ws.Range("H83").Value = 1476.25
$ws.Range("I83").Value = 1120.6666
$ws.Range("J83").Value = 2543
$ws.Range("K83").Value = 5603.333000000001
$ws.Range("L83").Value = 12715
$ws.Range("M83").Value = -611.3330000000005
$ws.Range("N83").Value = -22699
$ws.Range("H92").Value = 11666.286
$ws.Range("J92").Value = 11666.286
$ws.Range("L92").Value = 11666.286
$ws.Range("N92").Value = -15410.286
$ws.Range("H93").Value = 29900
$ws.Range("J93").Value = 29900
$ws.Range("L93").Value = 29900
$ws.Range("N93").Value = -33644
$ws.Range("H94").Value = 100000
$ws.Range("J94").Value = 100000
$ws.Range("L94").Value = 100000
$ws.Range("N94").Value = -101352
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("N95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("H97").Value = 7571.4287
$ws.Range("I97").Value = 10000
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -9504
$ws.Range("N97").Value = -2492
$ws.Range("H98").Value = 21599.4
$ws.Range("J98").Value = 21599.4
$ws.Range("L98").Value = 21599.4
$ws.Range("N98").Value = -27589.4
$ws.Range("H126").Value = 22225.562
$ws.Range("I126").Value = 26354.54
$ws.Range("K126").Value = 79063.62
$ws.Range("M126").Value = -76593.62
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").Value = 0
$ws.Range("L141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19555.223
$ws.Range("I7").Value = 22999.715
$ws.Range("K7").Value = 22999.715
$ws.Range("M7").Value = -22887.715
$ws.Range("H46").Value = 3712.125
$ws.Range("I46").Value = 1949.6666
$ws.Range("J46").Value = 8999.5
$ws.Range("K46").Value = 1949.6666
$ws.Range("L46").Value = 8999.5
$ws.Range("M46").Value = -1761.6666
$ws.Range("N46").Value = -9375.5
$ws.Range("H61").Value = 4704.2104
$ws.Range("I61").Value = 4092
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 4092
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -3890
$ws.Range("N61").Value = -7404
$ws.Range("H82").Value = 1065.7778
$ws.Range("I82").Value = 941.7143
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 941.7143
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -580.7143
$ws.Range("N82").Value = -2222
$ws.Range("H85").Value = 1065.7778
$ws.Range("I85").Value = 941.7143
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 941.7143
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = 306.2857
$ws.Range("N85").Value = -3996
$ws.Range("H93").Value = 465108
$ws.Range("I93").Value = 695954.7
$ws.Range("J93").Value = 3414.625
$ws.Range("K93").Value = 695954.7
$ws.Range("L93").Value = 3414.625
$ws.Range("M93").Value = -694706.7
$ws.Range("N93").Value = -5910.625
$ws.Range("H113").Value = 4704.2104
$ws.Range("I113").Value = 4092
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 4092
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = -1922
$ws.Range("N113").Value = -11340
$ws.Range("H126").Value = 19555.223
$ws.Range("I126").Value = 22999.715
$ws.Range("K126").Value = 68999.145
$ws.Range("M126").Value = -66529.145
$ws.Range("H136").Value = 6633.5713
$ws.Range("I136").Value = 6374.615
$ws.Range("K136").Value = 19123.845
$ws.Range("M136").Value = -16573.845
